$d = $word.ActiveDocument

# --- Simple 1:1 text replacements (date + most problem cells) ---
$replacements = @(
    @("2024-01-15 Monday", "2024-01-16 Tuesday"),
    @("240×3=", "790×5="),
    @("786×6=", "980×7="),
    @("896×2=", "153×5="),
    @("400×3=", "572×4="),
    @("820×7=", "709×4="),
    @("962×3=", "441×3="),
    @("762×5=", "905×6="),
    @("856×9=", "317×9="),
    @("636×7=", "560×4="),
    @("188×9=", "247×6="),
    @("211×3=", "910×8="),
    @("420×9=", "481×6="),
    @("364×2=", "779×4="),
    @("667×4=", "968×8="),
    @("383×7=", "706×4="),
    @("937×6=", "808×5="),
    @("125×4=", "175×8="),
    @("421×7=", "256×6="),
    @("863×4=", "882×9="),
    @("313×3=", "138×9=")
)

foreach ($pair in $replacements) {
    $d.Content.Find.Execute($pair[0], $true, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2) | Out-Null
}

# --- Structural change in row 1 of the table ---
# Old cells: 983×2=, 257×5=, 584×7=, 742×6=, 280×5=
# New cells: 152×5=, 685×6=, 141×4=, 120×6=, 742×6=
# The last old value (742×6=) is kept (now shifted one position to the
# right), a new cell (120×6=) is inserted, and the old last cell
# (280×5=) is dropped. Implement by inserting a brand-new row with the
# final content immediately before the old row, then deleting the old
# row.

$t = $d.Tables.Item(1)
$oldRow = $t.Rows.Item(1)
$newRow = $t.Rows.Add($oldRow)

$newValues = @("152×5=", "685×6=", "141×4=", "120×6=", "742×6=")
for ($i = 1; $i -le 5; $i++) {
    $newRow.Cells.Item($i).Range.Text = $newValues[$i - 1]
}

$t2 = $d.Tables.Item(1)
$t2.Rows.Item(2).Delete()
